# The "Interest details" section (old rows 61-75: Applicant interest,
# Owner details[], Interested persons[]) is removed entirely. This shifts
# the following sections (old rows 76-113) up by 15 rows, becoming the new
# rows 61-98, and a new "Interest in land" section header replaces the old
# "Interest details" text on what is now row 61-63 (previously
# applicant-owns-land / permission-obtained / permission-not-obtained-details).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old rows 61:75 (entire rows) - everything below shifts up.
$ws.Range("A61:N75").EntireRow.Delete() | Out-Null

# Update the new "Interest in land" section header (now rows 61-63, which
# used to be applicant-owns-land / permission-obtained /
# permission-not-obtained-details under the old "Interest details" section).
$ws.Range("A61").Value = "Interest in land"
$ws.Range("B61").Value = "interest-in-land"
$ws.Range("C61").Value = "Whether the applicant owns or has permission to use the land where the proposed advertisement will be placed"
$ws.Range("C62").Value = "Whether the applicant owns or has permission to use the land where the proposed advertisement will be placed"
$ws.Range("C63").Value = "Whether the applicant owns or has permission to use the land where the proposed advertisement will be placed"
